# Bestelformulier.xlsx update
# - Naam: changed from "Thomas van Bruchem" to "Bart van Wijk"
# - Klantennummer updated
# - Datum updated
# - Several order quantities (Aantal) cleared / updated (DB werkend, DB gekoppeld aan admin)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header fields
$ws.Range("B6").Value = "Bart van Wijk"
$ws.Range("B7").Value = 79001
$ws.Range("B8").Value = 44637

# Order quantity updates
$ws.Range("C11").Value = 2

$ws.Range("D17").ClearContents() | Out-Null
$ws.Range("D19").ClearContents() | Out-Null
$ws.Range("D33").ClearContents() | Out-Null
$ws.Range("C37").ClearContents() | Out-Null
$ws.Range("D39").ClearContents() | Out-Null
$ws.Range("C42").Value = 45
$ws.Range("C47").ClearContents() | Out-Null
$ws.Range("D55").ClearContents() | Out-Null
$ws.Range("D56").Value = 23
$ws.Range("D64").ClearContents() | Out-Null
$ws.Range("D65").ClearContents() | Out-Null
$ws.Range("C71").ClearContents() | Out-Null
$ws.Range("C76").Value = 45
$ws.Range("D86").ClearContents() | Out-Null
$ws.Range("D92").ClearContents() | Out-Null
$ws.Range("D94").ClearContents() | Out-Null
$ws.Range("D96").Value = 25
$ws.Range("D97").ClearContents() | Out-Null
$ws.Range("D108").Value = 4
$ws.Range("D109").ClearContents() | Out-Null
$ws.Range("C115").ClearContents() | Out-Null
$ws.Range("C117").ClearContents() | Out-Null
$ws.Range("C121").ClearContents() | Out-Null
$ws.Range("D124").ClearContents() | Out-Null
$ws.Range("D132").ClearContents() | Out-Null
$ws.Range("D138").Value = 23
$ws.Range("D139").ClearContents() | Out-Null
$ws.Range("D140").ClearContents() | Out-Null
$ws.Range("D146").ClearContents() | Out-Null
$ws.Range("C148").ClearContents() | Out-Null

# Restore selection to D17 as in the saved workbook
$ws.Range("D17").Select() | Out-Null

# Protect the worksheet (contents, objects and scenarios) as in the target workbook
$ws.Protect("", $true, $true, $true) | Out-Null
